# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows (Tuna - Vega Central Mapocho de Santiago)
# above the former last block of rows, pushing the existing rows 266-277
# down to 270-281.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 266; existing data shifts down to 270-281.
$ws.Rows("266:269").Insert()

# Columns A,B,C,E,F,G,H,I,K are constant across this product block.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100107
$producto  = "Otros"
$categoriaId = 100107011
$categoria = "Tuna"
$variedad  = "Sin especificar"
$unidad    = "`$/caja 18 kilos"
$origen    = "Provincia de Melipilla"

# New row 266: Calidad Especial
$r = 266
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44610
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 320
$ws.Cells.Item($r, 14).Value = 17000
$ws.Cells.Item($r, 15).Value = 17000
$ws.Cells.Item($r, 16).Value = 17000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 944
$ws.Cells.Item($r, 20).Value = 18

# New row 267: Calidad Primera
$r = 267
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44610
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 580
$ws.Cells.Item($r, 14).Value = 15000
$ws.Cells.Item($r, 15).Value = 15000
$ws.Cells.Item($r, 16).Value = 15000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 833
$ws.Cells.Item($r, 20).Value = 18

# New row 268: Calidad Segunda
$r = 268
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44610
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 400
$ws.Cells.Item($r, 14).Value = 12000
$ws.Cells.Item($r, 15).Value = 12000
$ws.Cells.Item($r, 16).Value = 12000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 667
$ws.Cells.Item($r, 20).Value = 18

# New row 269: Calidad Tercera
$r = 269
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44610
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Tercera"
$ws.Cells.Item($r, 13).Value = 260
$ws.Cells.Item($r, 14).Value = 8000
$ws.Cells.Item($r, 15).Value = 8000
$ws.Cells.Item($r, 16).Value = 8000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 444
$ws.Cells.Item($r, 20).Value = 18
